$d = $word.ActiveDocument

$d.Content.Find.Execute("726÷5=145, 1", $true, $false, $false, $false, $false, $true, 1, $false, "594÷9=66, 0", 2)
$d.Content.Find.Execute("839÷3=279, 2", $true, $false, $false, $false, $false, $true, 1, $false, "606÷9=67, 3", 2)
$d.Content.Find.Execute("903÷3=301, 0", $true, $false, $false, $false, $false, $true, 1, $false, "332÷7=47, 3", 2)
$d.Content.Find.Execute("735÷5=147, 0", $true, $false, $false, $false, $false, $true, 1, $false, "169÷4=42, 1", 2)
$d.Content.Find.Execute("716÷5=143, 1", $true, $false, $false, $false, $false, $true, 1, $false, "606÷2=303, 0", 2)
$d.Content.Find.Execute("430÷2=215, 0", $true, $false, $false, $false, $false, $true, 1, $false, "700÷2=350, 0", 2)
$d.Content.Find.Execute("408÷4=102, 0", $true, $false, $false, $false, $false, $true, 1, $false, "862÷3=287, 1", 2)
$d.Content.Find.Execute("695÷5=139, 0", $true, $false, $false, $false, $false, $true, 1, $false, "239÷9=26, 5", 2)
$d.Content.Find.Execute("982÷4=245, 2", $true, $false, $false, $false, $false, $true, 1, $false, "640÷5=128, 0", 2)
$d.Content.Find.Execute("540÷2=270, 0", $true, $false, $false, $false, $false, $true, 1, $false, "716÷4=179, 0", 2)
$d.Content.Find.Execute("980÷9=108, 8", $true, $false, $false, $false, $false, $true, 1, $false, "658÷5=131, 3", 2)
$d.Content.Find.Execute("357÷4=89, 1", $true, $false, $false, $false, $false, $true, 1, $false, "336÷3=112, 0", 2)
$d.Content.Find.Execute("620÷7=88, 4", $true, $false, $false, $false, $false, $true, 1, $false, "721÷2=360, 1", 2)
$d.Content.Find.Execute("410÷5=82, 0", $true, $false, $false, $false, $false, $true, 1, $false, "735÷2=367, 1", 2)
$d.Content.Find.Execute("968÷9=107, 5", $true, $false, $false, $false, $false, $true, 1, $false, "983÷3=327, 2", 2)
$d.Content.Find.Execute("974÷7=139, 1", $true, $false, $false, $false, $false, $true, 1, $false, "143÷7=20, 3", 2)
$d.Content.Find.Execute("857÷7=122, 3", $true, $false, $false, $false, $false, $true, 1, $false, "261÷9=29, 0", 2)
$d.Content.Find.Execute("568÷2=284, 0", $true, $false, $false, $false, $false, $true, 1, $false, "433÷2=216, 1", 2)
$d.Content.Find.Execute("515÷9=57, 2", $true, $false, $false, $false, $false, $true, 1, $false, "206÷3=68, 2", 2)
$d.Content.Find.Execute("723÷5=144, 3", $true, $false, $false, $false, $false, $true, 1, $false, "198÷3=66, 0", 2)
$d.Content.Find.Execute("908÷3=302, 2", $true, $false, $false, $false, $false, $true, 1, $false, "414÷7=59, 1", 2)
$d.Content.Find.Execute("130÷3=43, 1", $true, $false, $false, $false, $false, $true, 1, $false, "794÷7=113, 3", 2)
$d.Content.Find.Execute("733÷2=366, 1", $true, $false, $false, $false, $false, $true, 1, $false, "931÷8=116, 3", 2)
$d.Content.Find.Execute("877÷6=146, 1", $true, $false, $false, $false, $false, $true, 1, $false, "885÷3=295, 0", 2)
$d.Content.Find.Execute("299÷5=59, 4", $true, $false, $false, $false, $false, $true, 1, $false, "437÷2=218, 1", 2)
